$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Fix-Sql([string]$sql) {
    $sql = $sql -replace 'std\.id = prt\."study\.id"', 'std.study_id = prt."study.study_id"'
    $sql = $sql -replace 'prt\.id = dgn\."participant\.id"', 'prt.participant_id = dgn."participant.participant_id"'
    $sql = $sql -replace 'prt\.id = trt\."participant\.id"', 'prt.participant_id = trt."participant.participant_id"'
    $sql = $sql -replace 'prt\.id = trr\."participant\.id"', 'prt.participant_id = trr."participant.participant_id"'
    $sql = $sql -replace 'prt\.id = srv\."participant\.id"', 'prt.participant_id = srv."participant.participant_id"'
    $sql = $sql -replace 'std\.id = rfs\."study\.id"', 'std.study_id = rfs."study.study_id"'
    return $sql
}

$cells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $cell.Value = Fix-Sql $cell.Value()
}

$ws.Columns.Item(3).ColumnWidth = 68.15
